$d = $word.ActiveDocument

# --- Plain text / narrative edits -----------------------------------------

# 1. Timestamp in the header line "(09:07:3012 May 2020)" -> "(09:14:1512 May 2020)"
$d.Content.Find.Execute("09:07:3012 May 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "09:14:1512 May 2020", 2)

# 2. "estimated an odds ratio" -> "estimated a risk ratio"
$d.Content.Find.Execute("we estimated an odds ratio to", $true, $false, $false, $false, $false,
                         $true, 1, $false, "we estimated a risk ratio to", 2)

# 3. "logit link" -> "log link"
$d.Content.Find.Execute("binomial errors and logit link) to account", $true, $false, $false, $false, $false,
                         $true, 1, $false, "binomial errors and log link) to account", 2)

# --- Table of regression results -------------------------------------------
# Table 1, addressed as (row, column), 1-indexed, matching Word's Table.Cell API.

$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $oldValue, $newValue) {
    $cell = $table.Cell($row, $col)
    $current = $cell.Range.Text
    # Cell range text includes trailing cell-mark characters; compare the
    # visible text only.
    $visible = $current.TrimEnd([char]7, [char]13)
    if ($visible -ne $oldValue) {
        Write-Output ("WARNING: cell($row,$col) expected '" + $oldValue + "' but found '" + $visible + "'")
    }
    $cell.Range.Text = $newValue
}

# Row 4 ("D" / arm)
Set-CellValue $t 4 3 "0.06" "0.04"
Set-CellValue $t 4 4 "0.36" "0.37"
Set-CellValue $t 4 5 "0.72" "0.71"
Set-CellValue $t 4 6 "0.91" "0.93"
Set-CellValue $t 4 7 "1.15" "1.11"

# Row 7 ("8d9c30" / strat_var)
Set-CellValue $t 7 2 "1.14" "1.10"
Set-CellValue $t 7 3 "0.12" "0.09"
Set-CellValue $t 7 4 "1.26" "1.25"
Set-CellValue $t 7 6 "0.93" "0.94"
Set-CellValue $t 7 7 "1.39" "1.29"

# Row 8 ("9d5ed6" / strat_var)
Set-CellValue $t 8 2 "0.92" "0.94"
Set-CellValue $t 8 3 "0.12" "0.09"
Set-CellValue $t 8 6 "0.71" "0.77"
Set-CellValue $t 8 7 "1.19" "1.14"

# Row 9 ("e1e1d3" / strat_var)
Set-CellValue $t 9 2 "1.12" "1.09"
Set-CellValue $t 9 3 "0.12" "0.09"
Set-CellValue $t 9 4 "1.08" "1.07"
Set-CellValue $t 9 6 "0.91" "0.93"
Set-CellValue $t 9 7 "1.37" "1.27"

# Row 10 ("ff4457" / strat_var)
Set-CellValue $t 10 2 "1.21" "1.15"
Set-CellValue $t 10 3 "0.12" "0.09"
Set-CellValue $t 10 4 "1.84" "1.83"
Set-CellValue $t 10 7 "1.48" "1.34"

# Row 12 ("_cons")
Set-CellValue $t 12 2 "0.30" "0.23"
Set-CellValue $t 12 3 "0.03" "0.02"
Set-CellValue $t 12 4 "-13.28" "-21.19"
Set-CellValue $t 12 6 "0.25" "0.20"
Set-CellValue $t 12 7 "0.36" "0.27"

Write-Output "done"
